$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.575787544250488
$ws.Range("B1").Value = 4.401920795440674
$ws.Range("C1").Value = 3.10261344909668
$ws.Range("D1").Value = 1.247157335281372
$ws.Range("E1").Value = 0.9212923645973206
